$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-looking numeric cells remain text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.856.54"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "2.352.68"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "545.49"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").Value = "136.78"
$ws.Range("E6").Value = "  -3.25%  "

$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  -3.99%  "

$ws.Range("D9").Value = "2.353.21"
$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("E11").Value = "  +1.80%  "

$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("D13").Value = "0.346"
$ws.Range("E13").Value = "  +0.21%  "

$ws.Range("E14").Value = "  -3.26%  "

$ws.Range("D15").Value = "2.778.16"
$ws.Range("E15").Value = "  -1.09%  "

$ws.Range("D16").Value = "61.041.36"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("E17").Value = "  -1.17%  "

$ws.Range("D18").Value = "2.347.89"
$ws.Range("E18").Value = "  -1.32%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "320.35"
$ws.Range("E20").Value = "  +1.18%  "

$ws.Range("E21").Value = "  +0.66%  "

$ws.Range("D22").Value = "6.57"
$ws.Range("E22").Value = "  -2.10%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "63.42"
$ws.Range("E24").Value = "  +0.89%  "

$ws.Range("D25").Value = "1.67"
$ws.Range("E25").Value = "  -8.91%  "

$ws.Range("D26").Value = "8.49"
$ws.Range("E26").Value = "  +9.16%  "

$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.468.20"
$ws.Range("E28").Value = "  -1.16%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "8.01"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "1.38"
$ws.Range("E30").Value = "  -3.23%  "

$ws.Range("D31").Value = "0.0₃0869"
$ws.Range("E31").Value = "  -6.76%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "495.17"
$ws.Range("E32").Value = "  -5.03%  "

$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.146"
$ws.Range("E33").Value = "  +1.91%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "1.79"
$ws.Range("E34").Value = "  -2.19%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.51"
$ws.Range("E35").Value = "  -3.26%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "4.64"
$ws.Range("E37").Value = "  -0.55%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.377"
$ws.Range("E38").Value = "  +0.66%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "18.52"
$ws.Range("E39").Value = "  +2.88%  "

$ws.Range("E40").Value = "  -4.27%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  +5.83%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "145.02"
$ws.Range("E42").Value = "  +5.48%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "142.46"
$ws.Range("E44").Value = "  +1.75%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "3.58"
$ws.Range("E45").Value = "  +0.59%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  -7.40%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0516"
$ws.Range("E47").Value = "  -0.61%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "19.17"
$ws.Range("E48").Value = "  -5.81%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.569"
$ws.Range("E49").Value = "  -1.15%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.0902"
$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0221"
$ws.Range("E51").Value = "  -1.13%  "

# Restore default number format for the price column so no stray styling remains
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"